$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" timestamp in A1 (actually stored in A2) ---
$ws.Range("A2").Value = "Laatst bijgewerkt: 2025-09-06 22:44:13"

# --- 2. Update the advies-tabblad data cells that changed between runs ---
$ws.Range("H4").Value = "jacks"
$ws.Range("N4").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

$ws.Range("D5").Value = "meer dan 9.5"
$ws.Range("F5").Value = 3.5
$ws.Range("G5").Value = "minder dan 9.5"
$ws.Range("I5").Value = 1.44
$ws.Range("J5").Value = "1=44, 2=106"
$ws.Range("K5").Value = "€2.64"
$ws.Range("L5").Value = 1.98
$ws.Range("N5").Value = "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.44%7Creplace"

$ws.Range("H6").Value = "kambi"
$ws.Range("N6").Value = "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace"

$ws.Range("H7").Value = "jacks"
$ws.Range("N7").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

$ws.Range("D8").Value = "meer dan 8.5"
$ws.Range("F8").Value = 2.45
$ws.Range("G8").Value = "minder dan 8.5"
$ws.Range("I8").Value = 1.74
$ws.Range("J8").Value = "1=62, 2=88"
$ws.Range("K8").Value = "€1.9"
$ws.Range("L8").Value = 1.71
$ws.Range("N8").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace"

$ws.Range("H9").Value = "betmgm"
$ws.Range("N9").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace"

$ws.Range("D11").Value = "meer dan 7.5"
$ws.Range("F11").Value = 1.85
$ws.Range("G11").Value = "minder dan 7.5"
$ws.Range("H11").Value = "jacks"
$ws.Range("I11").Value = 2.25
$ws.Range("J11").Value = "1=82, 2=68"
$ws.Range("K11").Value = "€1.7"
$ws.Range("N11").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

$ws.Range("H12").Value = "betmgm"
$ws.Range("N12").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.43%7Creplace"

# --- 3. Rebuild the hyperlinks so M3:N12 point at the right targets/fragments ---
# Remember the "Hyperlink" cell style so it can be re-applied: Hyperlinks.Add()
# resets the style of the cell it is attached to.
$hlStyle = $ws.Range("M3").Style

$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("M3"), "https://sport.toto.nl/wedden/wedstrijd/8778584")
$ws.Hyperlinks.Add($ws.Range("N3"), "https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394")
$ws.Hyperlinks.Add($ws.Range("M4"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N4"), "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.44%7Creplace", "event/1023224945")
$ws.Hyperlinks.Add($ws.Range("M5"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N5"), "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace")
$ws.Hyperlinks.Add($ws.Range("M6"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N6"), "https://jacks.nl/sports/event/1023224945")
$ws.Hyperlinks.Add($ws.Range("M7"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N7"), "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace", "event/1023224945")
$ws.Hyperlinks.Add($ws.Range("M8"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N8"), "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace")
$ws.Hyperlinks.Add($ws.Range("M9"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N9"), "https://jacks.nl/sports/event/1023224945")
$ws.Hyperlinks.Add($ws.Range("M10"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N10"), "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace")
$ws.Hyperlinks.Add($ws.Range("M11"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N11"), "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.43%7Creplace", "event/1023224945")
$ws.Hyperlinks.Add($ws.Range("M12"), "https://sport.toto.nl/wedden/wedstrijd/8706282")
$ws.Hyperlinks.Add($ws.Range("N12"), "https://jacks.nl/sports/event/1023224945")

# Re-apply the original "Hyperlink" style and the display text, since adding
# the hyperlinks overwrote both with Excel's defaults.
$cols = @("M","N")
for ($r = 3; $r -le 12; $r++) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Style = $hlStyle
    }
}

$ws.Range("N4").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"
$ws.Range("N6").Value = "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace"
$ws.Range("N7").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"
$ws.Range("N9").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace"
$ws.Range("N11").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"
$ws.Range("N12").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.43%7Creplace"
